$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "The Journal of Heart and Lung Transplantation"
$ws.Range("G3").Value = "https://openalex.org/S42270594"
$ws.Range("H3").Value = "Elsevier BV"
$ws.Range("I3").Value = "1053-2498"
